# Auto-generated Excel COM-interop script applying scheduled-runner Sheets update
# Updates currentAveragePrice / LevePriceNQ/HQ / LeveProfitNQ/HQ columns (H,I,J,K,L,M,N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1053.75
$ws.Range("I53").Value = 905
$ws.Range("J53").Value = 1500
$ws.Range("K53").Value = 905
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = -268
$ws.Range("N53").Value = -2774
$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
$ws.Range("H116").Value = 9996.25
$ws.Range("J116").Value = 9996.25
$ws.Range("L116").Value = 9996.25
$ws.Range("N116").Value = -16880.25
$ws.Range("H132").Value = 5977.3
$ws.Range("I132").Value = 4295.5
$ws.Range("K132").Value = 12886.5
$ws.Range("M132").Value = -10356.5
$ws.Range("H137").Value = 841996
$ws.Range("J137").Value = 12082.833
$ws.Range("L137").Value = 36248.499
$ws.Range("N137").Value = -41348.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1684.5454
$ws.Range("I32").Value = 1684.5454
$ws.Range("K32").Value = 1684.5454
$ws.Range("M32").Value = -1397.5454
$ws.Range("H61").Value = 9499.6
$ws.Range("I61").Value = 4499.3335
$ws.Range("K61").Value = 4499.3335
$ws.Range("M61").Value = -4287.3335
$ws.Range("H74").Value = 2500
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2500
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 13292.143
$ws.Range("I132").Value = 11136.25
$ws.Range("K132").Value = 33408.75
$ws.Range("M132").Value = -30878.75
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 9499.6
$ws.Range("I136").Value = 4499.3335
$ws.Range("K136").Value = 13498.0005
$ws.Range("M136").Value = -10948.0005
$ws.Range("H137").Value = 72500
$ws.Range("J137").Value = 72500
$ws.Range("L137").Value = 72500
$ws.Range("N137").Value = -82700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1660
$ws.Range("I107").Value = 1575
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1575
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 345
$ws.Range("N107").Value = -5840
$ws.Range("H134").Value = 10019.75
$ws.Range("I134").Value = 5039.5
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 15118.5
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -12583.5
$ws.Range("N134").Value = -50070
$ws.Range("H135").Value = 52958.8
$ws.Range("J135").Value = 52958.8
$ws.Range("L135").Value = 52958.8
$ws.Range("N135").Value = -63098.8
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 53111.332
$ws.Range("J141").Value = 53111.332
$ws.Range("L141").Value = 53111.332
$ws.Range("N141").Value = -63471.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H59").Value = 34664.668
$ws.Range("I59").Value = 30497
$ws.Range("J59").Value = 43000
$ws.Range("K59").Value = 30497
$ws.Range("L59").Value = 43000
$ws.Range("M59").Value = -29352
$ws.Range("N59").Value = -45290
$ws.Range("H94").Value = 565.3333
$ws.Range("I94").Value = 348.5
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 348.5
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = 102.5
$ws.Range("N94").Value = -1901
$ws.Range("H132").Value = 9999.857
$ws.Range("J132").Value = 19333
$ws.Range("L132").Value = 57999
$ws.Range("N132").Value = -63059
$ws.Range("H134").Value = 6612.5454
$ws.Range("I134").Value = 1568.6
$ws.Range("K134").Value = 4705.799999999999
$ws.Range("M134").Value = -2170.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4478.4443
$ws.Range("I4").Value = 360.64285
$ws.Range("J4").Value = 8913
$ws.Range("K4").Value = 1081.92855
$ws.Range("L4").Value = 26739
$ws.Range("M4").Value = -969.9285500000001
$ws.Range("N4").Value = -26963
$ws.Range("H34").Value = 1200.6
$ws.Range("J34").Value = 1200.6
$ws.Range("L34").Value = 3601.8
$ws.Range("N34").Value = -3769.8
$ws.Range("H39").Value = 152
$ws.Range("J39").Value = 204
$ws.Range("L39").Value = 612
$ws.Range("N39").Value = -1200
$ws.Range("H92").Value = 1504.2
$ws.Range("I92").Value = 1750
$ws.Range("J92").Value = 1340.3334
$ws.Range("K92").Value = 5250
$ws.Range("L92").Value = 4021.0002
$ws.Range("M92").Value = -4002
$ws.Range("N92").Value = -6517.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 5000
$ws.Range("M102").Value = -3378
$ws.Range("H132").Value = 120109.11
$ws.Range("I132").Value = 132372.75
$ws.Range("J132").Value = 22000
$ws.Range("K132").Value = 397118.25
$ws.Range("L132").Value = 66000
$ws.Range("M132").Value = -394588.25
$ws.Range("N132").Value = -71060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9860
$ws.Range("N10").ClearContents()
$ws.Range("H16").Value = 1082.5
$ws.Range("J16").Value = 998.3333
$ws.Range("L16").Value = 998.3333
$ws.Range("N16").Value = -1338.3333
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 1700.5
$ws.Range("J61").Value = 2001
$ws.Range("L61").Value = 2001
$ws.Range("N61").Value = -2405
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 1700.5
$ws.Range("J113").Value = 2001
$ws.Range("L113").Value = 2001
$ws.Range("N113").Value = -6341
$ws.Range("H136").Value = 16999.8
$ws.Range("I136").Value = 12000
$ws.Range("K136").Value = 36000
$ws.Range("M136").Value = -33450
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 54000
$ws.Range("J129").Value = 54000
$ws.Range("L129").Value = 54000
$ws.Range("N129").Value = -64000
$ws.Range("H132").Value = 10999.875
$ws.Range("I132").Value = 8999.799999999999
$ws.Range("J132").Value = 14333.333
$ws.Range("K132").Value = 26999.4
$ws.Range("L132").Value = 42999.999
$ws.Range("M132").Value = -24469.4
$ws.Range("N132").Value = -48059.999
$ws.Range("H136").Value = 11600
$ws.Range("I136").Value = 9400
$ws.Range("J136").Value = 12333.333
$ws.Range("K136").Value = 28200
$ws.Range("L136").Value = 36999.999
$ws.Range("M136").Value = -25650
$ws.Range("N136").Value = -42099.999
